# Fruta / hortaliza, semanal
# Inserts two new weekly price rows into the Alcachofa (Agricola del Norte S.A. de Arica) table.
#
# Before the edit, rows 10-16 held (in order, newest-first by date):
#   10: 2021-08-04 Symphony
#   11: 2021-07-28 Madrigal
#   12: 2021-07-14 Madrigal
#   13: 2020-11-25 Madrigal
#   14: 2021-06-23 Argentina(o)
#   15: 2021-06-23 Madrigal
#   16: 2021-07-21 Madrigal
#
# After the edit, a new row is inserted at the very top (row 10: 2021-08-18 Madrigal)
# pushing the old rows 10-12 down to 11-13, then another new row is inserted
# (row 14: 2021-08-11 Symphony) pushing the remaining old rows 13-16 down to 15-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceRow {
    param($Row, $Fecha, $Variedad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $PrecioKg, $KgUnidades)

    $ws.Cells.Item($Row, 1).Value = 1
    $ws.Cells.Item($Row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($Row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value = 15
    $ws.Cells.Item($Row, 6).Value = 100112013
    $ws.Cells.Item($Row, 7).Value = "Alcachofa"
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = "Región de Coquimbo"
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# Insert the newest record (2021-08-18, Madrigal) at row 10, pushing everything else down
# one row (old row 10 -> 11, 11 -> 12, 12 -> 13, ...).
$ws.Rows(10).Insert()
Set-PriceRow 10 44426 "Madrigal" 150 19000 20000 19500 "`$/caja 40 unidades" 488 40

# Insert the second new record (2021-08-11, Symphony) at what is now row 14, pushing the
# remaining older rows (old 13-16, now at 13-16) down one more row to 15-18.
$ws.Rows(14).Insert()
Set-PriceRow 14 44419 "Symphony" 150 21000 22000 21500 "`$/caja 50 unidades" 430 50
